$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the (previously) blank separator row 88.
$ws.Rows("88:91").Insert()

# Row 87: the end time was corrected (0.91666... -> 0.92708...), which
# shifts the computed minutes/hours via the existing formulas.
$ws.Range("E87").Value = 0.92708333333333337

# New row 88: 2014-03-23
$ws.Range("A88").Value = 2014
$ws.Range("B88").Value = 3
$ws.Range("C88").Value = 23
$ws.Range("D88").Value = 0.57291666666666663
$ws.Range("E88").Value = 0.59375
$ws.Range("F88").Formula = "=(E88-D88)*24*60"
$ws.Range("G88").Formula = "=F88/60"

# New row 89: 2014-03-24
$ws.Range("A89").Value = 2014
$ws.Range("B89").Value = 3
$ws.Range("C89").Value = 24
$ws.Range("D89").Value = 0.3611111111111111
$ws.Range("E89").Value = 0.39583333333333331
$ws.Range("F89").Formula = "=(E89-D89)*24*60"
$ws.Range("G89").Formula = "=F89/60"

# New row 90: 2014-03-24
$ws.Range("A90").Value = 2014
$ws.Range("B90").Value = 3
$ws.Range("C90").Value = 24
$ws.Range("D90").Value = 0.71180555555555547
$ws.Range("E90").Value = 0.78125
$ws.Range("F90").Formula = "=(E90-D90)*24*60"
$ws.Range("G90").Formula = "=F90/60"

# New row 91: 2014-03-24
$ws.Range("A91").Value = 2014
$ws.Range("B91").Value = 3
$ws.Range("C91").Value = 24
$ws.Range("D91").Value = 0.82638888888888884
$ws.Range("E91").Value = 0.91666666666666663
$ws.Range("F91").Formula = "=(E91-D91)*24*60"
$ws.Range("G91").Formula = "=F91/60"

$ws.Range("A92").Select()
